$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "-"
$ws.Range("B6").Value = "-"
$ws.Range("B10").Value = "-"
$ws.Range("D10").Value = "-"
